# Fruta / hortaliza, semanal
# Insert a new weekly record at row 662 (pushing the existing rows 662-689
# down to 663-690) on the single data sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 662:689 down to 663:690, creating a blank (but formatted) row 662
$ws.Rows.Item(662).Insert()

# Populate the new row 662 with the new weekly price record
$ws.Cells.Item(662, 1).Value = 3
$ws.Cells.Item(662, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(662, 3).Value = "Coquimbo"
$ws.Cells.Item(662, 4).Value = 45147
$ws.Cells.Item(662, 5).Value = 5
$ws.Cells.Item(662, 6).Value = "Fruta"
$ws.Cells.Item(662, 7).Value = 100108
$ws.Cells.Item(662, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(662, 9).Value = 100108002
$ws.Cells.Item(662, 10).Value = "Mango"
$ws.Cells.Item(662, 11).Value = "Sin especificar"
$ws.Cells.Item(662, 12).Value = "Primera"
$ws.Cells.Item(662, 13).Value = 228
$ws.Cells.Item(662, 14).Value = 8000
$ws.Cells.Item(662, 15).Value = 8000
$ws.Cells.Item(662, 16).Value = 8000
$ws.Cells.Item(662, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(662, 18).Value = "Brasil"
$ws.Cells.Item(662, 19).Value = 2000
$ws.Cells.Item(662, 20).Value = 4
